# Update the cryptos price table (columns D = Price, E = Volume(1h)) with the
# latest scraped values. Column D entries are leading-apostrophe quoted so
# Excel stores them as literal text (matching the source data, where prices
# such as "26.595.60" / "1.943.00" are not valid numbers, and others such as
# "240.71" must still be kept as text rather than being auto-converted).
# In a single-quoted PowerShell string, '' is an escaped literal single
# quote, so '''26.595.60' evaluates to the text  '26.595.60 (a leading
# apostrophe followed by the price), which Excel interprets as "force text".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.595.60'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '''1.718.62'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''240.71'
$ws.Range('E5').Value = '  -2.18%  '
$ws.Range('D6').Value = '''0.9982'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '''0.4929'
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('D8').Value = '''0.2601'
$ws.Range('E8').Value = '  -2.80%  '
$ws.Range('D9').Value = '''0.06208'
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('D10').Value = '''1.726.62'
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('D11').Value = '''0.07003'
$ws.Range('D12').Value = '''15.76'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '''0.6073'
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('D14').Value = '''4.488'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').Value = '''76.83'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('D16').Value = '''0.9985'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = '''26.445.93'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '''0.9979'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '''0.000007155'
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('D21').Value = '''1.943.00'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').Value = '''4.415'
$ws.Range('E22').Value = '  -3.10%  '
$ws.Range('D23').Value = '''8.513'
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('D24').Value = '''5.086'
$ws.Range('E24').Value = '  -3.86%  '
$ws.Range('D25').Value = '''137.78'
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').Value = '''15.29'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('D27').Value = '''1.402'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').Value = '''105.80'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('D30').Value = '''3.919'
$ws.Range('E30').Value = '  -2.52%  '
$ws.Range('D31').Value = '''0.07952'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').Value = '''3.647'
$ws.Range('E32').Value = '  -2.05%  '
$ws.Range('D33').Value = '''0.04494'
$ws.Range('E33').Value = '  -2.56%  '
$ws.Range('D34').Value = '''2.611'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').Value = '''0.9992'
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('D37').Value = '''0.9363'
$ws.Range('E37').Value = '  +3.03%  '
$ws.Range('D38').Value = '''1.994'
$ws.Range('E38').Value = '  -2.72%  '
$ws.Range('D39').Value = '''2.414'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').Value = '''0.9980'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').Value = '''0.01510'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').Value = '''5.517'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = '''99.26'
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('E44').Value = '  -2.66%  '
$ws.Range('D45').Value = '''6.906'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('D47').Value = '''0.05375'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').Value = '''7.768'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').Value = '''30.04'
$ws.Range('E49').Value = '  -2.14%  '
$ws.Range('D50').Value = '''51.42'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('D51').Value = '''1.224'
$ws.Range('E51').Value = '  -2.31%  '
